$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Students Details "

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Grade"
$ws.Range("D1").Value = "Estimation"
$ws.Range("E1").Value = "Parent's Email"
$ws.Range("F1").Value = "Address"
$ws.Range("G1").Value = "Courses"
$ws.Range("H1").Value = "Age"
$ws.Range("I1").Value = "Phone"
$ws.Range("J1").Value = "isAbsent"

# Row 2 - Student1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Student1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "A"
$ws.Range("E2").Value = "parents_Student1@gmail.com"
$ws.Range("F2").Value = "1 Giza Egypt"
$ws.Range("G2").Value = "[English, Math]"
$ws.Range("H2").Value = 14
$ws.Range("I2").Value = 1111111111
$ws.Range("J2").Value = $true

# Row 3 - student2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "student2"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "4"
$ws.Range("D3").Value = "A"
$ws.Range("E3").Value = "parents_Strudent2@gmail.com"
$ws.Range("F3").Value = "2 Cairo"
$ws.Range("G3").Value = "[Arabic]"
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 222222222
$ws.Range("J3").Value = $true

# Row 4 - student5
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "student5"
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = "A"
$ws.Range("E4").Value = "parentStudent5@gmail"
$ws.Range("F4").Value = "55 street"
$ws.Range("G4").Value = "[Ar, En, Math, German]"
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 555555555
$ws.Range("J4").Value = $true

# Active cell selection per diff (D4)
$ws.Range("D4").Select()
